$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 0.134753
$ws.Cells.Item(2, 8).Value = 0.404259
$ws.Cells.Item(2, 9).Value = 0.005003875147349546
$ws.Cells.Item(2, 10).Value = 0.005003875147349547
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 2.885018666666667
$ws.Cells.Item(2, 14).Value = 8.655056
$ws.Cells.Item(2, 15).Value = 0.2739459408611671
$ws.Cells.Item(2, 16).Value = 0.2739459408611671
$ws.Cells.Item(2, 17).Value = 0.3887649203893333
$ws.Cells.Item(2, 18).Value = 3.498884283504
$ws.Cells.Item(2, 19).Value = 0.001370791285192483
$ws.Cells.Item(2, 20).Value = 0.001370791285192483
# Row 3
$ws.Cells.Item(3, 7).Value = 0.134753
$ws.Cells.Item(3, 8).Value = 0.404259
$ws.Cells.Item(3, 9).Value = 0.005003875147349546
$ws.Cells.Item(3, 10).Value = 0.005003875147349547
$ws.Cells.Item(3, 15).Value = 0.1360096554953469
$ws.Cells.Item(3, 16).Value = 0.1360096554953469
$ws.Cells.Item(3, 17).Value = 0.1930153910096666
$ws.Cells.Item(3, 18).Value = 1.737138519087
$ws.Cells.Item(3, 19).Value = 0.0006805753349327401
$ws.Cells.Item(3, 20).Value = 0.0006805753349327402
# Row 4
$ws.Cells.Item(4, 7).Value = 0.134753
$ws.Cells.Item(4, 8).Value = 0.404259
$ws.Cells.Item(4, 9).Value = 0.005003875147349546
$ws.Cells.Item(4, 10).Value = 0.005003875147349547
$ws.Cells.Item(4, 13).Value = 5.197721
$ws.Cells.Item(4, 14).Value = 15.593163
$ws.Cells.Item(4, 15).Value = 0.4935477839815871
$ws.Cells.Item(4, 16).Value = 0.4935477839815871
$ws.Cells.Item(4, 17).Value = 0.700408497913
$ws.Cells.Item(4, 18).Value = 6.303676481217
$ws.Cells.Item(4, 19).Value = 0.002469651490294906
$ws.Cells.Item(4, 20).Value = 0.002469651490294906
# Row 5
$ws.Cells.Item(5, 7).Value = 0.134753
$ws.Cells.Item(5, 8).Value = 0.404259
$ws.Cells.Item(5, 9).Value = 0.005003875147349546
$ws.Cells.Item(5, 10).Value = 0.005003875147349547
$ws.Cells.Item(5, 13).Value = 1.016239
$ws.Cells.Item(5, 14).Value = 3.048717
$ws.Cells.Item(5, 15).Value = 0.09649661966189875
$ws.Cells.Item(5, 16).Value = 0.09649661966189876
$ws.Cells.Item(5, 17).Value = 0.136941253967
$ws.Cells.Item(5, 18).Value = 1.232471285703
$ws.Cells.Item(5, 19).Value = 0.0004828570369294167
$ws.Cells.Item(5, 20).Value = 0.0004828570369294169
# Row 6
$ws.Cells.Item(6, 9).Value = 0.9088028687403782
$ws.Cells.Item(6, 10).Value = 0.9088028687403783
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 2.885018666666667
$ws.Cells.Item(6, 14).Value = 8.655056
$ws.Cells.Item(6, 15).Value = 0.2739459408611671
$ws.Cells.Item(6, 16).Value = 0.2739459408611671
$ws.Cells.Item(6, 17).Value = 70.6074121578738
$ws.Cells.Item(6, 18).Value = 635.466709420864
$ws.Cells.Item(6, 19).Value = 0.2489628569344107
$ws.Cells.Item(6, 20).Value = 0.2489628569344107
# Row 7
$ws.Cells.Item(7, 9).Value = 0.9088028687403782
$ws.Cells.Item(7, 10).Value = 0.9088028687403783
$ws.Cells.Item(7, 15).Value = 0.1360096554953469
$ws.Cells.Item(7, 16).Value = 0.1360096554953469
$ws.Cells.Item(7, 19).Value = 0.1236059650905618
$ws.Cells.Item(7, 20).Value = 0.1236059650905618
# Row 8
$ws.Cells.Item(8, 9).Value = 0.9088028687403782
$ws.Cells.Item(8, 10).Value = 0.9088028687403783
$ws.Cells.Item(8, 13).Value = 5.197721
$ws.Cells.Item(8, 14).Value = 15.593163
$ws.Cells.Item(8, 15).Value = 0.4935477839815871
$ws.Cells.Item(8, 16).Value = 0.4935477839815871
$ws.Cells.Item(8, 17).Value = 127.2080604430414
$ws.Cells.Item(8, 18).Value = 1144.872543987372
$ws.Cells.Item(8, 19).Value = 0.4485376419429228
$ws.Cells.Item(8, 20).Value = 0.4485376419429229
# Row 9
$ws.Cells.Item(9, 9).Value = 0.9088028687403782
$ws.Cells.Item(9, 10).Value = 0.9088028687403783
$ws.Cells.Item(9, 13).Value = 1.016239
$ws.Cells.Item(9, 14).Value = 3.048717
$ws.Cells.Item(9, 15).Value = 0.09649661966189875
$ws.Cells.Item(9, 16).Value = 0.09649661966189876
$ws.Cells.Item(9, 17).Value = 24.87124494303867
$ws.Cells.Item(9, 18).Value = 223.841204487348
$ws.Cells.Item(9, 19).Value = 0.08769640477248276
$ws.Cells.Item(9, 20).Value = 0.08769640477248279
# Row 10
$ws.Cells.Item(10, 7).Value = 2.321161
$ws.Cells.Item(10, 8).Value = 6.963483
$ws.Cells.Item(10, 9).Value = 0.08619325611227224
$ws.Cells.Item(10, 10).Value = 0.08619325611227226
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 2.885018666666667
$ws.Cells.Item(10, 14).Value = 8.655056
$ws.Cells.Item(10, 15).Value = 0.2739459408611671
$ws.Cells.Item(10, 16).Value = 0.2739459408611671
$ws.Cells.Item(10, 17).Value = 6.696592813338667
$ws.Cells.Item(10, 18).Value = 60.269335320048
$ws.Cells.Item(10, 19).Value = 0.02361229264156397
$ws.Cells.Item(10, 20).Value = 0.02361229264156397
# Row 11
$ws.Cells.Item(11, 7).Value = 2.321161
$ws.Cells.Item(11, 8).Value = 6.963483
$ws.Cells.Item(11, 9).Value = 0.08619325611227224
$ws.Cells.Item(11, 10).Value = 0.08619325611227226
$ws.Cells.Item(11, 15).Value = 0.1360096554953469
$ws.Cells.Item(11, 16).Value = 0.1360096554953469
$ws.Cells.Item(11, 17).Value = 3.324748228324333
$ws.Cells.Item(11, 18).Value = 29.922734054919
$ws.Cells.Item(11, 19).Value = 0.01172311506985235
$ws.Cells.Item(11, 20).Value = 0.01172311506985235
# Row 12
$ws.Cells.Item(12, 7).Value = 2.321161
$ws.Cells.Item(12, 8).Value = 6.963483
$ws.Cells.Item(12, 9).Value = 0.08619325611227224
$ws.Cells.Item(12, 10).Value = 0.08619325611227226
$ws.Cells.Item(12, 13).Value = 5.197721
$ws.Cells.Item(12, 14).Value = 15.593163
$ws.Cells.Item(12, 15).Value = 0.4935477839815871
$ws.Cells.Item(12, 16).Value = 0.4935477839815871
$ws.Cells.Item(12, 17).Value = 12.064747274081
$ws.Cells.Item(12, 18).Value = 108.582725466729
$ws.Cells.Item(12, 19).Value = 0.04254049054836935
$ws.Cells.Item(12, 20).Value = 0.04254049054836936
# Row 13
$ws.Cells.Item(13, 7).Value = 2.321161
$ws.Cells.Item(13, 8).Value = 6.963483
$ws.Cells.Item(13, 9).Value = 0.08619325611227224
$ws.Cells.Item(13, 10).Value = 0.08619325611227226
$ws.Cells.Item(13, 13).Value = 1.016239
$ws.Cells.Item(13, 14).Value = 3.048717
$ws.Cells.Item(13, 15).Value = 0.09649661966189875
$ws.Cells.Item(13, 16).Value = 0.09649661966189876
$ws.Cells.Item(13, 17).Value = 2.358854333479
$ws.Cells.Item(13, 18).Value = 21.229689001311
$ws.Cells.Item(13, 19).Value = 0.008317357852486565
$ws.Cells.Item(13, 20).Value = 0.008317357852486567
